$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('KHACH_HANG')
$ws2 = $wb.Worksheets.Item('UY_QUYEN')

# KHACH_HANG row 9
$ws1.Cells.Item(9,1).Value = 8
$ws1.Cells.Item(9,2).Value = 'HẢI'
$ws1.Cells.Item(9,3).Value = '222222222222222222222222222222222222222222233333333333333333333'
$ws1.Cells.Item(9,4).Value = '333@gmail.com'
$ws1.Cells.Item(9,5).Value = '1111111111111111111111111111111111111111222222222'
$ws1.Cells.Item(9,6).Value = ''''
$ws1.Cells.Item(9,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(9,8).Value = ''''
$ws1.Cells.Item(9,9).Value = '3333333333333333333333333333333333333333333333333344444444444444444'
$ws1.Cells.Item(9,10).Value = ''''
$ws1.Cells.Item(9,11).Value = 'index 10'

# KHACH_HANG row 10
$ws1.Cells.Item(10,1).Value = 9
$ws1.Cells.Item(10,2).Value = 'PHAN ANH QUÂN'
$ws1.Cells.Item(10,3).Value = '0912345678'
$ws1.Cells.Item(10,4).Value = 'quan.ap@gmail.com'
$ws1.Cells.Item(10,5).Value = '100010234'
$ws1.Cells.Item(10,6).Value = '21/10/2021'
$ws1.Cells.Item(10,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(10,8).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws1.Cells.Item(10,9).Value = '210123578'
$ws1.Cells.Item(10,10).Value = 'Techcombank'
$ws1.Cells.Item(10,11).Value = 'index 10'
$ws1.Cells.Item(10,12).Value = 'VP P30'
$ws1.Cells.Item(10,13).Value = 'test device'

# KHACH_HANG row 11
$ws1.Cells.Item(11,1).Value = 10
$ws1.Cells.Item(11,2).Value = 'NGUYỄN THỊ HOÀNG YẾN'
$ws1.Cells.Item(11,3).Value = '0934678359'
$ws1.Cells.Item(11,4).Value = 'yen.nguyenhoang@gmail.com'
$ws1.Cells.Item(11,5).Value = '210123321'
$ws1.Cells.Item(11,6).Value = '17/07/2020'
$ws1.Cells.Item(11,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(11,8).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws1.Cells.Item(11,9).Value = '0866730311'
$ws1.Cells.Item(11,10).Value = 'Vietcombank'
$ws1.Cells.Item(11,11).Value = 'index 10'
$ws1.Cells.Item(11,12).Value = 'Exim D210'
$ws1.Cells.Item(11,13).Value = 'VP P30'
$ws1.Cells.Item(11,14).Value = 'test device'

# KHACH_HANG row 12
$ws1.Cells.Item(12,1).Value = 11
$ws1.Cells.Item(12,2).Value = 'NGUYỄN THỊ HOÀNG YẾN'
$ws1.Cells.Item(12,3).Value = '0934678359'
$ws1.Cells.Item(12,4).Value = 'yen.nguyenhoang@gmail.com'
$ws1.Cells.Item(12,5).Value = '210123321'
$ws1.Cells.Item(12,6).Value = '17/07/2020'
$ws1.Cells.Item(12,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(12,8).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws1.Cells.Item(12,9).Value = '0866730311'
$ws1.Cells.Item(12,10).Value = 'Vietcombank'
$ws1.Cells.Item(12,11).Value = 'index 10'
$ws1.Cells.Item(12,12).Value = 'VP P30'

# KHACH_HANG row 13
$ws1.Cells.Item(13,1).Value = 12
$ws1.Cells.Item(13,2).Value = 'NGUYỄN THỊ HOÀNG YẾN'
$ws1.Cells.Item(13,3).Value = '0934678359'
$ws1.Cells.Item(13,4).Value = 'yen.nguyenhoang@gmail.com'
$ws1.Cells.Item(13,5).Value = '210123321'
$ws1.Cells.Item(13,6).Value = '17/07/2020'
$ws1.Cells.Item(13,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(13,8).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws1.Cells.Item(13,9).Value = '0866730311'
$ws1.Cells.Item(13,10).Value = 'Vietcombank'
$ws1.Cells.Item(13,11).Value = 'index 10'
$ws1.Cells.Item(13,12).Value = 'VP P30'
$ws1.Cells.Item(13,13).Value = 'Exim D210'

# KHACH_HANG row 14
$ws1.Cells.Item(14,1).Value = 13
$ws1.Cells.Item(14,2).Value = 'NGUYỄN THỊ HOÀNG YẾN'
$ws1.Cells.Item(14,3).Value = '0934678359'
$ws1.Cells.Item(14,4).Value = 'yen.nguyenhoang@gmail.com'
$ws1.Cells.Item(14,5).Value = '210123321'
$ws1.Cells.Item(14,6).Value = '17/07/2020'
$ws1.Cells.Item(14,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(14,8).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws1.Cells.Item(14,9).Value = '0866730311'
$ws1.Cells.Item(14,10).Value = 'Vietcombank'
$ws1.Cells.Item(14,11).Value = 'index 10'
$ws1.Cells.Item(14,12).Value = 'VP P30'
$ws1.Cells.Item(14,13).Value = 'Exim D210'
$ws1.Cells.Item(14,14).Value = 'Exim D210'

# KHACH_HANG row 15
$ws1.Cells.Item(15,1).Value = 14
$ws1.Cells.Item(15,2).Value = 'NGUYỄN THỊ HOÀNG YẾN'
$ws1.Cells.Item(15,3).Value = '0934678359'
$ws1.Cells.Item(15,4).Value = 'yen.nguyenhoang@gmail.com'
$ws1.Cells.Item(15,5).Value = '210123321'
$ws1.Cells.Item(15,6).Value = '17/07/2020'
$ws1.Cells.Item(15,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(15,8).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws1.Cells.Item(15,9).Value = '0866730311'
$ws1.Cells.Item(15,10).Value = 'Vietcombank'
$ws1.Cells.Item(15,11).Value = 'index 10'
$ws1.Cells.Item(15,12).Value = 'Exim D210'

# KHACH_HANG row 16
$ws1.Cells.Item(16,1).Value = 15
$ws1.Cells.Item(16,2).Value = 'NGUYỄN THỊ HOÀNG YẾN'
$ws1.Cells.Item(16,3).Value = '0934678359'
$ws1.Cells.Item(16,4).Value = 'yen.nguyenhoang@gmail.com'
$ws1.Cells.Item(16,5).Value = '210123321'
$ws1.Cells.Item(16,6).Value = '17/07/2020'
$ws1.Cells.Item(16,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(16,8).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws1.Cells.Item(16,9).Value = '0866730311'
$ws1.Cells.Item(16,10).Value = 'Vietcombank'
$ws1.Cells.Item(16,11).Value = 'index 10'
$ws1.Cells.Item(16,12).Value = 'Exim D210'

# KHACH_HANG row 17
$ws1.Cells.Item(17,1).Value = 16
$ws1.Cells.Item(17,2).Value = 'NGUYỄN THỊ HOÀNG YẾN'
$ws1.Cells.Item(17,3).Value = '0934678359'
$ws1.Cells.Item(17,4).Value = 'yen.nguyenhoang@gmail.com'
$ws1.Cells.Item(17,5).Value = '210123321'
$ws1.Cells.Item(17,6).Value = '17/07/2020'
$ws1.Cells.Item(17,7).Value = 'Cục cảnh sát QLHC về TTXH'
$ws1.Cells.Item(17,8).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws1.Cells.Item(17,9).Value = '0866730311'
$ws1.Cells.Item(17,10).Value = 'Vietcombank'
$ws1.Cells.Item(17,11).Value = 'index 10'
$ws1.Cells.Item(17,12).Value = 'VP P30'

# UY_QUYEN row 4
$ws2.Cells.Item(4,1).Value = 3
$ws2.Cells.Item(4,2).Value = 'CÔNG TY TNHH MTV PHAN THỊ'
$ws2.Cells.Item(4,3).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws2.Cells.Item(4,4).Value = '02363847G'
$ws2.Cells.Item(4,5).Value = '13/10/2020'
$ws2.Cells.Item(4,6).Value = 'Phòng Tài chính Kế hoạch Quận Hải Châu'
$ws2.Cells.Item(4,7).Value = 'PHAN HOÀNG HẢI'
$ws2.Cells.Item(4,8).Value = '111 Nguyễn Du, Phường Thạch Thang, Quận Hải Châu, Thành phố Đà Nẵng'
$ws2.Cells.Item(4,9).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws2.Cells.Item(4,10).Value = '0201585869'
$ws2.Cells.Item(4,11).Value = '13/07/2010'
$ws2.Cells.Item(4,12).Value = 'Công an Đà Nẵng'

# UY_QUYEN row 5
$ws2.Cells.Item(5,1).Value = 4
$ws2.Cells.Item(5,2).Value = ''''
$ws2.Cells.Item(5,3).Value = ''''
$ws2.Cells.Item(5,4).Value = 'ádấd'
$ws2.Cells.Item(5,5).Value = ''''
$ws2.Cells.Item(5,6).Value = 'áđấsđá'
$ws2.Cells.Item(5,7).Value = ''''
$ws2.Cells.Item(5,8).Value = 'áđasađâsd'
$ws2.Cells.Item(5,9).Value = ''''
$ws2.Cells.Item(5,10).Value = ''''
$ws2.Cells.Item(5,11).Value = ''''
$ws2.Cells.Item(5,12).Value = ''''

# UY_QUYEN row 6
$ws2.Cells.Item(6,1).Value = 5
$ws2.Cells.Item(6,2).Value = 'ZXÁD'
$ws2.Cells.Item(6,3).Value = 'áđâs'
$ws2.Cells.Item(6,4).Value = 'áđấ'
$ws2.Cells.Item(6,5).Value = 'áđâsd'
$ws2.Cells.Item(6,6).Value = ''''
$ws2.Cells.Item(6,7).Value = ''''
$ws2.Cells.Item(6,8).Value = ''''
$ws2.Cells.Item(6,9).Value = ''''
$ws2.Cells.Item(6,10).Value = ''''
$ws2.Cells.Item(6,11).Value = ''''
$ws2.Cells.Item(6,12).Value = ''''

# UY_QUYEN row 7
$ws2.Cells.Item(7,1).Value = 6
$ws2.Cells.Item(7,2).Value = ''''
$ws2.Cells.Item(7,3).Value = ''''
$ws2.Cells.Item(7,4).Value = 'áđâsd'
$ws2.Cells.Item(7,5).Value = ''''
$ws2.Cells.Item(7,6).Value = 'áđâsd'
$ws2.Cells.Item(7,7).Value = ''''
$ws2.Cells.Item(7,8).Value = ''''
$ws2.Cells.Item(7,9).Value = ''''
$ws2.Cells.Item(7,10).Value = ''''
$ws2.Cells.Item(7,11).Value = ''''
$ws2.Cells.Item(7,12).Value = ''''

# UY_QUYEN row 8
$ws2.Cells.Item(8,1).Value = 7
$ws2.Cells.Item(8,2).Value = 'TẬP ĐOÀN XÂY DỰNG VÀ ĐẦU TƯ PHAN THỊ'
$ws2.Cells.Item(8,3).Value = '111 Nguyễn Du, Phường Thạch Thang, Quận Hải Châu, Thành phố Đà Nẵng'
$ws2.Cells.Item(8,4).Value = '041173206H'
$ws2.Cells.Item(8,5).Value = '13/02/2010'
$ws2.Cells.Item(8,6).Value = 'Phòng Tài chính Kế hoạch Ủy ban nhân dân Quận Hải Châu'
$ws2.Cells.Item(8,7).Value = 'PHAN HOÀNG HẢI'
$ws2.Cells.Item(8,8).Value = '111 Nguyễn Du, Phường Thạch Thang, Quận Hải Châu, Thành phố Đà Nẵng'
$ws2.Cells.Item(8,9).Value = 'K7/6 Xuân Tâm, Phường Thuận Phước, Quận Hải Châu, Thành phố Đà Nẵng'
$ws2.Cells.Item(8,10).Value = '210585867'
$ws2.Cells.Item(8,11).Value = '13/10/2000'
$ws2.Cells.Item(8,12).Value = 'Công an Đà Nẵng'

